# Apply updated dSF (column F) values per repulled data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 2
    6  = 2
    11 = 2
    14 = 0
    19 = 2
    21 = 0
    22 = 1
    34 = 3
    36 = 3
    46 = -6
    67 = -2
    68 = 1
    69 = 0
    77 = 0
    81 = 3
    85 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
